# Update EC (Estado de Cuenta) database - "Actualiza base de datos EC y agrega
# parte 1 de nuevos estado de cuenta"
#
# The old worker list (MEREDITH LIZETH PERALTA IREGUI / JULIO CESAR JIMENEZ
# MENDOZA / ROBERTO CARLOS DUARTE GONZALEZ / RAFAEL GUSTAVO CARMONA ALONSO,
# 2 periods each = 8 rows) is replaced with a smaller, updated list (only
# ROBERTO CARLOS DUARTE GONZALEZ and RAFAEL GUSTAVO CARMONA ALONSO, 2 periods
# each = 4 rows), and the summary header cells are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give the soon-to-be-last data row (19) the heavier "closing" bottom
#     border that the old last row (23) had, before that row disappears. ---
$ws.Range("B23:J23").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122) # xlPasteFormats

# --- Remove the now-obsolete extra data rows (20-23); rows below (the
#     signature block, formerly 28/29) shift up automatically to 24/25. ---
$ws.Range("A20:A23").EntireRow.Delete()

# --- Refresh the remaining 4 data rows (16-19) with the updated roster:
#     ROBERTO CARLOS DUARTE GONZALEZ and RAFAEL GUSTAVO CARMONA ALONSO,
#     each for periods 1603 and 1604. ---
$ws.Range("C16").Value2 = "9297629"
$ws.Range("D16").Value2 = "ROBERTO CARLOS DUARTE GONZALEZ"
$ws.Range("E16").Value2 = "1603"
$ws.Range("F16").Value2 = 27600
$ws.Range("G16").Value2 = 689455

$ws.Range("C17").Value2 = "1047439971"
$ws.Range("D17").Value2 = "RAFAEL GUSTAVO CARMONA ALONSO"
$ws.Range("E17").Value2 = "1603"
$ws.Range("F17").Value2 = 27600
$ws.Range("G17").Value2 = 689455

$ws.Range("C18").Value2 = "9297629"
$ws.Range("D18").Value2 = "ROBERTO CARLOS DUARTE GONZALEZ"
$ws.Range("E18").Value2 = "1604"
$ws.Range("F18").Value2 = 27600
$ws.Range("G18").Value2 = 689455

$ws.Range("C19").Value2 = "1047439971"
$ws.Range("D19").Value2 = "RAFAEL GUSTAVO CARMONA ALONSO"
$ws.Range("E19").Value2 = "1604"
$ws.Range("F19").Value2 = 27600
$ws.Range("G19").Value2 = 689455

# --- Refresh the summary totals: total mora value and worker count. ---
$ws.Range("E11").Value2 = 110400
$ws.Range("C13").Value2 = 2
